# Fix extra time sum: the "مجموع" (total) row on the "گزارش تردد" sheet
# was sitting two rows below the header (row 4) with an empty row in
# between (row 2/3 had no data). Remove the blank row so the totals row
# moves up to row 3 directly beneath the data, eliminating the stray gap
# that was throwing off the extra-time sum layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("گزارش تردد")

# Row 2 is completely empty; deleting it shifts the "مجموع" row from 4 to 3.
$ws.Rows.Item(2).Delete()
